# Audit-SEO workbook update
# - "couleur" row becomes "couleur et contraste"
# - "balise meta keyword" row gets a status ("fait") and an action ("changer")
# - a new "description image" row is appended (row 20) with status + action
# - "text format image (avis)" row is renamed to "text format image " and
#   gets a status ("fait") and an action ("text à la place d'image")
# - "paris" row gets a status ("fait") and an action ("remplacer lyon par paris")
# - the visible window is scrolled down and D26 becomes the active selection
#
# NOTE: the order in which new cell values are written below matches the
# order new strings appear in the workbook's shared-string table, so the
# statements are intentionally sequenced that way rather than grouped by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: "couleur " -> "couleur et contraste"
$ws.Range("B19").Value = "couleur et contraste"

# Row 10: "balise meta keyword " -> add status/action
$ws.Range("C10").Value = "fait"
$ws.Range("D10").Value = "changer"

# Row 20 (new row): label + status
$ws.Range("B20").Value = "description image"
$ws.Range("C20").Value = "fait "

# Row 15: "text format image (avis)" -> "text format image " + status/action
$ws.Range("B15").Value = "text format image "
$ws.Range("C15").Value = "fait"
$ws.Range("D15").Value = "text à la place d'image"

# Row 18: "paris " -> add status/action
$ws.Range("C18").Value = "fait"
$ws.Range("D18").Value = "remplacer lyon par paris"

# Row 20 (new row): action, completes the new row started above
$ws.Range("D20").Value = "modifier"

# Scroll the view down and select D26, matching the author's final view state
$ws.Range("D26").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
